$wb = $excel.ActiveWorkbook

# "Revisão 2" is the sheet the new "Revisão 3" should be placed after.
$ws2 = $wb.Worksheets.Item("Revisão 2")

# Clone "Revisão 2" (same column widths / styles / page margins) and drop the
# clone right after it, then rename it. Copy() makes the clone the active sheet.
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$new = $wb.ActiveSheet
$new.Name = "Revisão 3"

# Overwrite the cloned data with the Revisão 3 rows (8 requirements instead of
# the 7 that "Revisão 2" had).
$labels = @("RFUN5.1","RFUN5.2","RFUN5.3","RFUN5.4","RFUN5.5","RFUN6.1","RFUN6.2","RFUN6.3")
$rowVals = @("Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $new.Cells.Item($row, 1).Value = $labels[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 2
        $new.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}

# Row 8 is new (the template only had 7 rows) so it needs the same centered
# style ("s=1" in the template) explicitly applied.
$lastRow = $new.Range("A8:K8")
$lastRow.HorizontalAlignment = -4108
$lastRow.VerticalAlignment = -4108

$new.Range("A9").Select()

Write-Output "done"
